$d = $word.ActiveDocument

# --- Change 1 & 3: the "_GoBack" bookmark moves from its old spot near --
# --- "增廣見聞" to a new spot wrapping the "啊」" run. Re-adding a ------
# --- bookmark with a name that already exists simply relocates it, so ---
# --- this single Add both removes the old bookmark and creates the ------
# --- new one. ------------------------------------------------------------
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("啊」", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $bmRange1 = $d.Range($rng1.Start, $rng1.End)
    $d.Bookmarks.Add("_GoBack", $bmRange1)
}

# --- Change 2: remove the stand-alone paragraph that contains only a ----
# --- page break (the one right before "伍、其他") ------------------------
# Search for the tail of the preceding paragraph, its paragraph mark, the
# page-break character and the page-break paragraph's own mark, so the
# match spans exactly "...得心應手。<PARA MARK><PAGE BREAK><PARA MARK>".
$anchorText = "中更得心應手。"
$pattern2 = $anchorText + [char]13 + [char]12 + [char]13
$rng2 = $d.Content
$found2 = $rng2.Find.Execute($pattern2, $false, $false, $false, $false, $false, $false, 1, $false, "", 0)
if ($found2) {
    # Keep the anchor text + its paragraph mark untouched; delete only the
    # trailing "<PAGE BREAK><PARA MARK>" pair, i.e. the empty paragraph.
    $delStart = $rng2.Start + $anchorText.Length + 1
    $delEnd = $rng2.End
    $delRange = $d.Range($delStart, $delEnd)
    $delRange.Delete()
}
